# Adds a new "LCD cursor" keyboard-shortcut entry to the keyboard-shortcuts
# sheet: a new row for key "u" -> "Toggle LCD cursor (off/profile/template)"
# inserted right after the "z" / "Toggle xy cursor clamp mode" row (row 11),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11 (shifts rows 11:28 down to 12:29,
# carrying their formatting/row-heights/styles down with them).
$ws.Rows("11:11").Insert()

# Fill in the new shortcut entry.
$ws.Cells.Item(11, 1).Value = "u"
$ws.Cells.Item(11, 2).Value = "Toggle LCD cursor (off/profile/template)"

# The newly inserted row picks up the (slightly condensed) row height used
# by its neighbouring rows instead of the sheet default.
$ws.Rows("11:11").RowHeight = 13.8

# Move the active selection to the new row, matching the author's cursor
# position after typing the new entry.
$ws.Range("B11").Select()
